$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a blank separator row above the "add Time..." open task ---
# (old row 48 "...add Time..." shifts down to become row 49, etc.)
$ws.Rows.Item(48).Insert()

# --- 2) Insert 9 rows before "In Fatigue module:" to make room for the new
#        "In Mechanical testing module:" section (rows 66-74 after step 1) ---
$ws.Rows("66:74").Insert()

# --- 3) Fill in the new "In Mechanical testing module:" open-task block ---
# (cell order matters for shared-string table layout, so fill it the same
#  way it was authored: A69, A70, B70, A71, A72, B71, A73, B73, B72)
$ws.Range("A69").Value = "In Mechanical testing module:"

$ws.Range("A70").Value = "CantileverBending and TwoPointBending"
$ws.Range("B70").Value = "are subclasses of BendingMeasurement -> but aren't they methods, and for the measurement process itself it should be, e.g., CantileverBendingMeasurement?"

$ws.Range("A71").Value = "MechanicalTesting"
$ws.Range("A72").Value = "TensileTesting"

$ws.Range("B71").Value = "same as CantileverBending"

$ws.Range("A73").Value = "TensileTestMeasurement"
$ws.Range("B73").Value = 'exists, but what is the difference to TensileTesting, and why are they both a "Measurement"?'

$ws.Range("B72").Value = "same as CantileverBending"

# --- 4) Update the selection to match the final state ---
$ws.Range("B74").Select()
